$d = $word.ActiveDocument

# --- 1. Locate the paragraph that ends with "This will search for all cities that end with s." ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "This will search for all cities that end with s.`r") {
        $targetIdx = $i
        break
    }
}
if ($targetIdx -eq -1) {
    throw "Could not find anchor paragraph"
}

# --- 2. Find an existing paragraph that already uses the bullet list (numId=2), so we can
#        reuse its list template for our new "SELECT * FROM Name" bullet item. ---
$srcListPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -eq "SELECT *`r" -and $pp.Range.ListFormat.ListType -eq 2) {
        $srcListPara = $pp
        break
    }
}
if ($null -eq $srcListPara) {
    throw "Could not find source list paragraph"
}
$tmpl = $srcListPara.Range.ListFormat.ListTemplate

$cur = $targetIdx

# --- 3. Two blank paragraphs ---
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

# --- 4. "IN" (bold) + "- is used to specify multiple values while using WHERE" ---
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "IN- is used to specify multiple values while using WHERE"
$p = $d.Paragraphs($cur)
$boldRange = $d.Range($p.Range.Start, $p.Range.Start + 2)
$boldRange.Font.Bold = 1

# --- 5. "SELECT * FROM Name" bulleted list item (numId=2) ---
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.Text = "SELECT * FROM Name"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($tmpl, $true, 2, $false, 1)

# --- 6. "WHERE City IN (‘fresno’, ‘visalia’);" continuation list paragraph ---
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.Text = "WHERE City IN (‘fresno’, ‘visalia’);"

# --- 7. Trailing empty ListParagraph-styled paragraph ---
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "List Paragraph"

Write-Host "Done. Total paragraphs now: $($d.Paragraphs.Count)"
